$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly completed progress-report row (row 3)
$ws.Range("A3").Value = Get-Date -Year 2017 -Month 1 -Day 30 -Hour 0 -Minute 0 -Second 0
$ws.Range("A3").NumberFormat = "m/d/yyyy"

$ws.Range("B3").Value = 4

$ws.Range("C3").Value = "Requirement Analysis Drafting"
